$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The filer added a new "latest period" column right after column C (the
# labels column). Insert a new column D, which shifts the previous D:K
# data/dates one column to the right (becoming E:L). Scope the insert to the
# table's row range so the sheet's used range/dimension doesn't balloon out
# to the full column.
$ws.Range("D5:D102").EntireColumn.Insert()

# The newly inserted column D does not inherit the number formatting of the
# data table (it instead picks up General format). Copy the formats (only)
# from the now-adjacent column E, which holds what used to be column D, so
# the new column D matches the rest of the table (date format on the
# "Period Ending" rows, #,##0 on the data rows). Only the rows that actually
# held data (7:102) are copied, so blank label-only rows (5, 6) don't pick
# up a spurious formatted-but-empty D cell.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 37 and 79 are bare section-title rows ("Balance Sheet", "Cash Flow
# Statement") that only ever had a label in column B and no D:K cells at
# all. The format-only paste above still stamped an empty, styled D cell on
# them (because the destination range spans every row from 7 to 102); clear
# those two back out so the row doesn't gain a cell it never had.
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# Populate the new column D with the latest reporting period's figures.
$colD = @{
    7 = 43465
    8 = 92200
    9 = "NA"
    10 = "NA"
    12 = "NA"
    13 = 0
    14 = 0
    15 = -1000
    17 = 18500
    18 = 73700
    20 = -52400
    21 = "NA"
    22 = 0
    23 = 21300
    24 = 3200
    25 = 0
    26 = 18100
    27 = 18100
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 52400
    33 = 18100
    34 = 0
    35 = 18100
    38 = 43465
    41 = 40000
    42 = 87300
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 56000
    49 = 79000
    50 = 0
    51 = 0
    52 = 5600
    53 = 0
    54 = 2274400
    57 = "NA"
    58 = 0
    59 = 7300
    60 = 0
    61 = 39200
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 1991400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 40000
    73 = 0
    74 = 0
    75 = 0
    76 = 283000
    77 = 0
    80 = 43465
    81 = 18100
    83 = "NA"
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 20800
    91 = -3800
    92 = 0
    93 = 0
    94 = -94900
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 76800
    101 = 0
    102 = 2800
}

foreach ($row in $colD.Keys) {
    $ws.Cells.Item($row, 4).Value = $colD[$row]
}

# Match column D's width to the rest of the data columns (D:K) now that it
# holds real figures, instead of leaving it at the sheet's default width.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth
